$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.847.94"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.641.83"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'216.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").Value = "'19.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").Value = "'0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.867.70"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "1.640.49"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "'0.562"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'63.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "25.879.84"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'4.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").Value = "'193.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "'9.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("E24").Value = "  +4.72%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'142.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "'0.0495"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "'3.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").Value = "1.132.23"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.548"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").Value = "'100.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D45").Value = "1.776.80"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("E46").Value = "  +3.06%  "
$ws.Range("D47").Value = "'55.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E49").Value = "  +6.03%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'2.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.49%  "
